$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.056.33"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "3.679.61"
$ws.Range("E3").Value = "  +7.37%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "583.65"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "177.23"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "3.671.27"
$ws.Range("E7").Value = "  +7.36%  "
$ws.Range("D8").Value = "0.619"
$ws.Range("E8").Value = "  +3.88%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "0.199"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "6.90"
$ws.Range("E11").Value = "  +26.72%  "
$ws.Range("D12").Value = "0.610"
$ws.Range("E12").Value = "  +4.53%  "
$ws.Range("D13").Value = "48.94"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "0.0000288"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "4.282.37"
$ws.Range("E15").Value = "  +7.77%  "
$ws.Range("D16").Value = "674.13"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").Value = "9.01"
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("D18").Value = "3.693.32"
$ws.Range("E18").Value = "  +7.65%  "
$ws.Range("D19").Value = "71.224.06"
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "17.92"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").Value = "11.53"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").Value = "0.939"
$ws.Range("E23").Value = "  +4.57%  "
$ws.Range("D24").Value = "17.29"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").Value = "101.73"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "3.97"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("D27").Value = "2.82"
$ws.Range("E27").Value = "  +6.09%  "
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  +5.27%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "34.99"
$ws.Range("E30").Value = "  +4.46%  "
$ws.Range("D31").Value = "3.40"
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("D32").Value = "9.10"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").Value = "7.56"
$ws.Range("E34").Value = "  +6.29%  "
$ws.Range("D35").Value = "4.03"
$ws.Range("E35").Value = "  +7.99%  "
$ws.Range("D36").Value = "584.68"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").Value = "11.19"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").Value = "0.108"
$ws.Range("E38").Value = "  +4.58%  "
$ws.Range("D39").Value = "58.71"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "0.0462"
$ws.Range("E41").Value = "  +10.07%  "
$ws.Range("D42").Value = "3.618.92"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "0.143"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").Value = "0.350"
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("D45").Value = "0.0₃0759"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("D46").Value = "35.10"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "2.75"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("E48").Value = "  +9.10%  "
$ws.Range("E49").Value = "  +3.90%  "
$ws.Range("D50").Value = "134.74"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").Value = "2.96"
$ws.Range("E51").Value = "  +8.17%  "
